$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.870.85"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.636.74"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.58"
$ws.Range("E5").Value = "  -2.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.06"
$ws.Range("E6").Value = "  -3.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.631.33"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  -4.97%  "
$ws.Range("E11").Value = "  +16.27%  "
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.59"
$ws.Range("E13").Value = "  -4.10%  "
$ws.Range("E14").Value = "  -2.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.222.72"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "675.67"
$ws.Range("E16").Value = "  -4.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.97"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.635.18"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.854.89"
$ws.Range("E19").Value = "  -2.29%  "
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("E21").Value = "  -4.22%  "
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.942"
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.17"
$ws.Range("E24").Value = "  -4.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "100.17"
$ws.Range("E25").Value = "  -4.89%  "
$ws.Range("E26").Value = "  -2.76%  "
$ws.Range("E27").Value = "  -2.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.83"
$ws.Range("E29").Value = "  -2.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.71"
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("E32").Value = "  -5.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.58"
$ws.Range("E33").Value = "  +1.84%  "
$ws.Range("E34").Value = "  -6.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.98"
$ws.Range("E35").Value = "  -4.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "576.51"
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("E37").Value = "  -2.38%  "
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "58.46"
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.567.17"
$ws.Range("E41").Value = "  -2.12%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0453"
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.346"
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("E44").Value = "  -3.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "34.38"
$ws.Range("E45").Value = "  -4.47%  "
$ws.Range("E46").Value = "  -6.21%  "
$ws.Range("E47").Value = "  -4.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.95"
$ws.Range("E48").Value = "  +4.20%  "
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.73"
$ws.Range("E50").Value = "  +2.93%  "
$ws.Range("E51").Value = "  -4.16%  "
